$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two address values that changed
$ws.Range("I2").Value = "2993 Highway 100, Villa Ridge, MO"
$ws.Range("Q2").Value = "1419 Highway Am, Villa Ridge, MO 63089"

# Autofit column I to reflect the new (shorter) content width
$ws.Columns.Item(9).ColumnWidth = 28.8

# Update the view: scroll so column E is the top-left visible column, and change selection to F3
$ws.Range("F3").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 5
